# Automation of valid login tc
# - Rename the worksheet from "tc1" to "ValidLogin"
# - Move the active selection from D16 to A2

$wb = $excel.ActiveWorkbook

# Grab the first worksheet (works whether or not it has already been renamed)
$ws = $wb.Worksheets.Item(1)

$ws.Name = "ValidLogin"

$ws.Activate()
$ws.Range("A2").Select()
